$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.464.27'
$ws.Range('E2').Value = '  -2.57%  '
$ws.Range('D3').Value = '1.988.87'
$ws.Range('E3').Value = '  -1.24%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''238.02'
$ws.Range('E5').Value = '  -9.43%  '
$ws.Range('E6').Value = '  -3.28%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''54.67'
$ws.Range('E8').Value = '  -2.58%  '
$ws.Range('D9').Value = '''0.371'
$ws.Range('E9').Value = '  -3.94%  '
$ws.Range('D10').Value = '''58.47'
$ws.Range('E10').Value = '  +3.21%  '
$ws.Range('D11').Value = '''0.0749'
$ws.Range('E11').Value = '  -3.35%  '
$ws.Range('D12').Value = '''0.0987'
$ws.Range('E12').Value = '  -2.99%  '
$ws.Range('D13').Value = '''14.23'
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('D14').Value = '2.280.31'
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('D15').Value = '''20.06'
$ws.Range('E15').Value = '  -4.75%  '
$ws.Range('D16').Value = '''0.754'
$ws.Range('E16').Value = '  -6.58%  '
$ws.Range('D17').Value = '''5.06'
$ws.Range('E17').Value = '  -3.85%  '
$ws.Range('D18').Value = '1.989.29'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D19').Value = '36.446.71'
$ws.Range('E19').Value = '  -2.25%  '
$ws.Range('D20').Value = '''67.81'
$ws.Range('E20').Value = '  -2.78%  '
$ws.Range('D21').Value = '0.0₃0804'
$ws.Range('E21').Value = '  -4.74%  '
$ws.Range('D22').Value = '''5.28'
$ws.Range('E22').Value = '  +1.91%  '
$ws.Range('D23').Value = '''221.48'
$ws.Range('E23').Value = '  -3.05%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '''2.41'
$ws.Range('E25').Value = '  -10.45%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '''2.37'
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('D27').Value = '''162.93'
$ws.Range('E27').Value = '  -1.09%  '
$ws.Range('D28').Value = '''8.70'
$ws.Range('E28').Value = '  -3.44%  '
$ws.Range('D29').Value = '''0.129'
$ws.Range('E29').Value = '  +0.46%  '
$ws.Range('D30').Value = '''18.88'
$ws.Range('E30').Value = '  -4.18%  '
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').Value = '''0.116'
$ws.Range('E32').Value = '  -2.97%  '
$ws.Range('D33').Value = '''4.36'
$ws.Range('E33').Value = '  -6.34%  '
$ws.Range('D34').Value = '''0.0607'
$ws.Range('E34').Value = '  -6.57%  '
$ws.Range('D35').Value = '''4.24'
$ws.Range('E35').Value = '  -6.53%  '
$ws.Range('E36').Value = '  -2.51%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('E39').Value = '  -2.65%  '
$ws.Range('D40').Value = '''5.46'
$ws.Range('E40').Value = '  +5.04%  '
$ws.Range('D41').Value = '''3.00'
$ws.Range('E41').Value = '  -1.56%  '
$ws.Range('D42').Value = '1.461.48'
$ws.Range('E42').Value = '  +4.71%  '
$ws.Range('D43').Value = '''0.0921'
$ws.Range('E43').Value = '  -2.48%  '
$ws.Range('E44').Value = '  -5.31%  '
$ws.Range('E45').Value = '  -10.09%  '
$ws.Range('D46').Value = '''88.91'
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('D47').Value = '''14.92'
$ws.Range('E47').Value = '  -5.23%  '
$ws.Range('D48').Value = '''0.994'
$ws.Range('E48').Value = '  -3.18%  '
$ws.Range('D49').Value = '''2.88'
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('D50').Value = '''6.79'
$ws.Range('E50').Value = '  -4.12%  '
$ws.Range('E51').Value = '  +16.69%  '
